# Despesas com consorcios updated
# Append 16 new rows (366-381) to the ConsorcioDespesas table on the
# "Despesas" sheet, mirroring the same column layout/formulas used by the
# existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# consorcio | data_base | funcao | subfuncao | ndo | empenhado | liquidado | pago
$newData = @(
    @("COFRON", 45626, 4,  122, 319011010100, 1355.05, 1355.05, 1355.05),
    @("COFRON", 45626, 4,  122, 319013010100, 72.11,   72.11,   72.03),
    @("COFRON", 45626, 4,  122, 319013020100, 189.29,  189.29,  189.08),
    @("COFRON", 45626, 4,  122, 339014140000, 0,       0,       0),
    @("COFRON", 45626, 4,  122, 339030000000, 64.07,   64.07,   64.07),
    @("COFRON", 45626, 4,  122, 339033010000, 0,       0,       0),
    @("COFRON", 45626, 4,  122, 339039000000, 31.47,   89.63,   81.8),
    @("COFRON", 45626, 4,  122, 339039990100, 59.11,   59.11,   59.11),
    @("COFRON", 45626, 4,  122, 339039400000, 247.52,  111.6,   111.6),
    @("COFRON", 45626, 4,  122, 339046010100, 53.48,   53.48,   53.48),
    @("COFRON", 45626, 4,  122, 339047000000, 0,       0,       0),
    @("COFRON", 45626, 4,  122, 339049010000, 1.96,    1.96,    1.96),
    @("COFRON", 45626, 4,  122, 449052000000, 0,       0,       0),
    @("COFRON", 45626, 10, 302, 334041390500, 9273.68, 9273.68, 9273.68),
    @("COFRON", 45626, 10, 302, 334041391100, 865.54,  865.54,  865.54),
    @("COFRON", 45626, 10, 302, 334041391000, 655.91,  655.91,  655.91)
)

$lastRow = $tbl.Range.Row + $tbl.Range.Rows.Count - 1
$firstNewRow = 0
$lastNewRow = 0

for ($i = 0; $i -lt $newData.Count; $i++) {
    $row = $newData[$i]
    $listRow = $tbl.ListRows.Add()
    $r = $listRow.Range
    $rowNum = $r.Row

    if ($i -eq 0) { $firstNewRow = $rowNum }
    $lastNewRow = $rowNum

    # Copy number formats from the row directly above so the new row reuses
    # the same style indices (date / ndo custom format / accounting format)
    # instead of Excel minting brand-new ones.
    $srcRow = $ws.Range("A" + ($rowNum - 1) + ":K" + ($rowNum - 1))
    $srcRow.Copy()
    $r.PasteSpecial(-4122)
    $excel.CutCopyMode = $false

    $r.Cells.Item(1, 1).Value = $row[0]
    $r.Cells.Item(1, 2).Value = $row[1]
    $r.Cells.Item(1, 3).Value = $row[2]
    $r.Cells.Item(1, 4).Value = $row[3]
    $r.Cells.Item(1, 5).Value = $row[4]
    $r.Cells.Item(1, 6).Value = $row[5]
    $r.Cells.Item(1, 7).Value = $row[6]
    $r.Cells.Item(1, 8).Value = $row[7]

    $r.Cells.Item(1, 9).Formula = "=YEAR(ConsorcioDespesas[[#This Row],[data_base]])"
    # Table (ListObject) calculated columns reject array-entered formulas
    # (same restriction Excel itself enforces), so this is entered as a
    # normal formula; it still evaluates to the same cached result as the
    # single-cell dynamic-array form used elsewhere in the column.
    $r.Cells.Item(1, 10).Formula = "=_xlfn.SWITCH(MONTH(ConsorcioDespesas[[#This Row],[data_base]]),1,1,2,1,3,2,4,2,5,3,6,3,7,4,8,4,9,5,10,5,11,6,12,6)"
    $r.Cells.Item(1, 11).Formula = "=MONTH(ConsorcioDespesas[[#This Row],[data_base]])"
}

$wb.Application.CutCopyMode = $false

# Leave the selection/scroll position where the author left off editing.
$ws.Range("F379").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 346
$win.ScrollColumn = 1
